# sprint6EstimatedTimes.xlsx - "Updating spreadsheet with new Task 4952"
#
# A new task row (Task Id 4952, "Fix display issues on Front- End UI") is
# inserted into the "Agilefant Timesheet" sheet as the last task under the
# "U35 - Activity search - (Modified ACs)" story (row 33), pushing every
# row below it down by one. The TOTALS formulas at the bottom of the sheet
# automatically extend to include the new row once it is inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33 (existing rows 33+ shift down to 34+)
$ws.Rows.Item(33).Insert()

# Fill in the new row with the new task's data
$ws.Range("A33").Value = 422
$ws.Range("B33").Value = "U35 - Activity search - (Modified ACs)"
$ws.Range("C33").Value = 4952
$ws.Range("D33").Value = "Fix display issues on Front- End UI"
$ws.Range("E33").Value = 2

# Leave the selection on the newly-entered row, matching the saved cursor
# position left behind by the edit.
$ws.Range("F33").Select() | Out-Null
